$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. The existing "总计" sheet (sheetId 6) becomes the new "2022-Q1" detail
#    sheet - rename it in place so it keeps its identity/position.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Use the "2021-Q4" sheet (same fund-holding-detail layout) as a formatting
# template so headers/borders/fonts match the rest of the workbook exactly.
$template = $wb.Worksheets.Item("2021-Q4")

$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$template.Range("A2:H9").Copy()
$q1.Range("A2:H9").PasteSpecial(-4122)

# ---- header row ----
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B-G on the data rows are stored as number-looking TEXT strings
# (matching every other quarterly sheet in this workbook). Assign them as
# ="..." formulas, then convert the whole block to plain values in one shot,
# so the result is text with no numeric coercion and no stray
# NumberFormat-driven style index left behind.

$q1.Range("A2").Value = 0
$q1.Range("B2").Formula = "=""008969"""
$q1.Range("C2").Formula = "=""睿远均衡价值三年持有期混合A"""
$q1.Range("D2").Formula = "=""141.91"""
$q1.Range("E2").Formula = "=""93.01"""
$q1.Range("F2").Formula = "=""3.49"""
$q1.Range("G2").Formula = "=""4.9527"""
$q1.Range("H2").Value = 7

$q1.Range("A3").Value = 1
$q1.Range("B3").Formula = "=""008970"""
$q1.Range("C3").Formula = "=""睿远均衡价值三年持有期混合C"""
$q1.Range("D3").Formula = "=""15.26"""
$q1.Range("E3").Formula = "=""93.01"""
$q1.Range("F3").Formula = "=""3.49"""
$q1.Range("G3").Formula = "=""0.5326"""
$q1.Range("H3").Value = 7

$q1.Range("A4").Value = 2
$q1.Range("B4").Formula = "=""006973"""
$q1.Range("C4").Formula = "=""太平睿盈混合A"""
$q1.Range("D4").Formula = "=""8.69"""
$q1.Range("E4").Formula = "=""29.33"""
$q1.Range("F4").Formula = "=""2.33"""
$q1.Range("G4").Formula = "=""0.2025"""
$q1.Range("H4").Value = 1

$q1.Range("A5").Value = 3
$q1.Range("B5").Formula = "=""010268"""
$q1.Range("C5").Formula = "=""太平睿安混合A"""
$q1.Range("D5").Formula = "=""5.52"""
$q1.Range("E5").Formula = "=""40.12"""
$q1.Range("F5").Formula = "=""2.97"""
$q1.Range("G5").Formula = "=""0.1639"""
$q1.Range("H5").Value = 1

$q1.Range("A6").Value = 4
$q1.Range("B6").Formula = "=""007669"""
$q1.Range("C6").Formula = "=""太平睿盈混合C"""
$q1.Range("D6").Formula = "=""2.21"""
$q1.Range("E6").Formula = "=""29.33"""
$q1.Range("F6").Formula = "=""2.33"""
$q1.Range("G6").Formula = "=""0.0515"""
$q1.Range("H6").Value = 1

$q1.Range("A7").Value = 5
$q1.Range("B7").Formula = "=""014053"""
$q1.Range("C7").Formula = "=""太平睿庆混合A"""
$q1.Range("D7").Formula = "=""2.46"""
$q1.Range("E7").Formula = "=""31.00"""
$q1.Range("F7").Formula = "=""1.21"""
$q1.Range("G7").Formula = "=""0.0298"""
$q1.Range("H7").Value = 8

$q1.Range("A8").Value = 6
$q1.Range("B8").Formula = "=""014054"""
$q1.Range("C8").Formula = "=""太平睿庆混合C"""
$q1.Range("D8").Formula = "=""0.67"""
$q1.Range("E8").Formula = "=""31.00"""
$q1.Range("F8").Formula = "=""1.21"""
$q1.Range("G8").Formula = "=""0.0081"""
$q1.Range("H8").Value = 8

$q1.Range("A9").Value = 7
$q1.Range("B9").Formula = "=""010269"""
$q1.Range("C9").Formula = "=""太平睿安混合C"""
$q1.Range("D9").Formula = "=""0.07"""
$q1.Range("E9").Formula = "=""40.12"""
$q1.Range("F9").Formula = "=""2.97"""
$q1.Range("G9").Formula = "=""0.0021"""
$q1.Range("H9").Value = 1

$q1.Range("B2:G9").Copy()
$q1.Range("B2:G9").PasteSpecial(-4163)

# ---------------------------------------------------------------------------
# 2. Create a brand-new "总计" sheet at the end of the workbook holding the
#    quarter-over-quarter summary (old data plus the new 2022-Q1 row). Copy the
#    whole "2021-Q4" sheet (rather than Worksheets.Add) so sheetPr/pageMargins
#    match the conventions already used by the other sheets, then wipe it.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy([System.Reflection.Missing]::Value, $lastSheet)
$total = $wb.Worksheets.Item($wb.Worksheets.Count)
$total.Name = "总计"
$total.Cells.Clear()

# Re-use the same header / column-A formatting used elsewhere in the workbook.
$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$q1.Range("A2:D7").Copy()
$total.Range("A2:D7").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Formula = "=""2022-Q1"""
$total.Range("C2").Value = 8
$total.Range("D2").Value = 5.94

$total.Range("A3").Value = 1
$total.Range("B3").Formula = "=""2021-Q4"""
$total.Range("C3").Value = 8
$total.Range("D3").Value = 6.21

$total.Range("A4").Value = 2
$total.Range("B4").Formula = "=""2021-Q3"""
$total.Range("C4").Value = 8
$total.Range("D4").Value = 8.5

$total.Range("A5").Value = 3
$total.Range("B5").Formula = "=""2021-Q2"""
$total.Range("C5").Value = 7
$total.Range("D5").Value = 13.27

$total.Range("A6").Value = 4
$total.Range("B6").Formula = "=""2021-Q1"""
$total.Range("C6").Value = 17
$total.Range("D6").Value = 15.66

$total.Range("A7").Value = 5
$total.Range("B7").Formula = "=""2020-Q4"""
$total.Range("C7").Value = 17
$total.Range("D7").Value = 28.12

$total.Range("B2:B7").Copy()
$total.Range("B2:B7").PasteSpecial(-4163)

